# Handback status report: add a new row for file
# 7332775b-08c9-4371-bbb2-94e10f68287a.md (in sync with en-US) to every
# worksheet: Overview, zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$fileGuid   = "7332775b-08c9-4371-bbb2-94e10f68287a"
$fileName   = "$fileGuid.md"
$pathName   = "e2e\$fileGuid.md"
$statusSync = "Handed back: in sync with en-US"

# -------------------------------------------------------------------
# Sheet "Overview" (columns A..G) - processed strictly left to right
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30877432d1026706d7e805da846a32c3bb81e3c2/$pathName",
    "",
    "",
    $pathName
) | Out-Null

$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G4").Value = "2016-08-21 14:52:06"

# -------------------------------------------------------------------
# Sheet "zh-cn" (columns A..P) - processed strictly left to right
# -------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$zhXlf = "$fileGuid.3acf4f0a054ca90b834b0c4c31fcc519b47cdd02.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30877432d1026706d7e805da846a32c3bb81e3c2/$pathName",
    "",
    "",
    $fileName
) | Out-Null

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H4").Value = "2016-08-21 14:51:58"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e48478dcb74f21345d2cce8038a39d5e0853964b/$pathName",
    "",
    "",
    $fileName
) | Out-Null

$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K4").Value = "2016-08-21 14:52:27"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

# -------------------------------------------------------------------
# Sheet "de-de" (columns A..P) - processed strictly left to right
# -------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$deXlf = "$fileGuid.3acf4f0a054ca90b834b0c4c31fcc519b47cdd02.de-de.xlf"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30877432d1026706d7e805da846a32c3bb81e3c2/$pathName",
    "",
    "",
    $fileName
) | Out-Null

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H4").Value = "2016-08-21 14:52:06"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/21a278f64f7fd633dbdde131ca3766e4d58e72e3/$pathName",
    "",
    "",
    $fileName
) | Out-Null

$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K4").Value = "2016-08-21 14:52:33"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wb.Save()
